$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)

# --- Insert a new row 2 into the summary (总计) sheet for 2022-Q3 ---
$sheetTotal.Rows.Item(2).Insert()
# Copy formatting (style) from the row below, which still holds the old row-2 look
$sheetTotal.Cells.Item(3,1).Copy($sheetTotal.Cells.Item(2,1))

# Fill in the new 2022-Q3 summary values
$sheetTotal.Cells.Item(2,1).Value = 0
$sheetTotal.Cells.Item(2,2).Value = "2022-Q3"
$sheetTotal.Cells.Item(2,3).Value = 12
$sheetTotal.Cells.Item(2,4).Value = 2

# Renumber the sequence column (A) for every pushed-down data row
for ($r = 3; $r -le 9; $r++) {
    $sheetTotal.Cells.Item($r,1).Value = $r - 2
}

# --- Create the new 2022-Q3 worksheet (copy 2022-Q2 to inherit its styling) ---
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($null, $sheetTotal)
$ws = $wb.ActiveSheet
$ws.Name = "2022-Q3"

# 2022-Q3 only lists 12 funds (2022-Q2 had 13); drop the extra inherited row
$ws.Rows.Item(14).Delete()

# Overwrite data rows with the 2022-Q3 figures
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "512980"
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "广发中证传媒ETF"
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "44.76"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "99.29"
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = "3.10"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "1.3876"
$ws.Cells.Item(2,8).Value = 8
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,2).Value = "160629"
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "鹏华中证传媒指数（LOF）A"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "6.41"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "94.58"
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = "2.92"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "0.1872"
$ws.Cells.Item(3,8).Value = 8
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = "002270"
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "东吴安盈量化灵活配置混合A"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "4.41"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "52.91"
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = "2.51"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "0.1107"
$ws.Cells.Item(4,8).Value = 8
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,2).Value = "516620"
$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,3).Value = "国泰中证影视主题ETF"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.94"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "99.07"
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = "9.06"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "0.0852"
$ws.Cells.Item(5,8).Value = 2
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).NumberFormat = "@"
$ws.Cells.Item(6,2).Value = "159855"
$ws.Cells.Item(6,3).NumberFormat = "@"
$ws.Cells.Item(6,3).Value = "银华中证影视主题ETF"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "0.84"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "96.84"
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = "8.84"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "0.0743"
$ws.Cells.Item(6,8).Value = 2
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).NumberFormat = "@"
$ws.Cells.Item(7,2).Value = "159805"
$ws.Cells.Item(7,3).NumberFormat = "@"
$ws.Cells.Item(7,3).Value = "鹏华中证传媒ETF"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "1.71"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "98.37"
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = "3.06"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "0.0523"
$ws.Cells.Item(7,8).Value = 8
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).NumberFormat = "@"
$ws.Cells.Item(8,2).Value = "164818"
$ws.Cells.Item(8,3).NumberFormat = "@"
$ws.Cells.Item(8,3).Value = "工银瑞信中证传媒指数（LOF）A"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "1.65"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "93.46"
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = "2.89"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "0.0477"
$ws.Cells.Item(8,8).Value = 8
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).NumberFormat = "@"
$ws.Cells.Item(9,2).Value = "290012"
$ws.Cells.Item(9,3).NumberFormat = "@"
$ws.Cells.Item(9,3).Value = "泰信行业精选灵活配置混合A"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.75"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "91.96"
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = "5.27"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "0.0395"
$ws.Cells.Item(9,8).Value = 9
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).NumberFormat = "@"
$ws.Cells.Item(10,2).Value = "010677"
$ws.Cells.Item(10,3).NumberFormat = "@"
$ws.Cells.Item(10,3).Value = "工银瑞信中证传媒指数（LOF）C"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.21"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "93.46"
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = "2.89"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "0.0061"
$ws.Cells.Item(10,8).Value = 8
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).NumberFormat = "@"
$ws.Cells.Item(11,2).Value = "015675"
$ws.Cells.Item(11,3).NumberFormat = "@"
$ws.Cells.Item(11,3).Value = "鹏华中证传媒指数（LOF）C"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.17"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "94.58"
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = "2.92"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = "0.0050"
$ws.Cells.Item(11,8).Value = 8
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).NumberFormat = "@"
$ws.Cells.Item(12,2).Value = "015154"
$ws.Cells.Item(12,3).NumberFormat = "@"
$ws.Cells.Item(12,3).Value = "东吴安盈量化灵活配置混合C"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.19"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "52.91"
$ws.Cells.Item(12,6).NumberFormat = "@"
$ws.Cells.Item(12,6).Value = "2.51"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = "0.0048"
$ws.Cells.Item(12,8).Value = 8
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).NumberFormat = "@"
$ws.Cells.Item(13,2).Value = "002583"
$ws.Cells.Item(13,3).NumberFormat = "@"
$ws.Cells.Item(13,3).Value = "泰信行业精选灵活配置混合C"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.04"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "91.96"
$ws.Cells.Item(13,6).NumberFormat = "@"
$ws.Cells.Item(13,6).Value = "5.27"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = "0.0021"
$ws.Cells.Item(13,8).Value = 9

# Restore the original active sheet/selection (总计)
$sheetTotal.Activate()
[void]$sheetTotal.Range("A1").Select()

Write-Output "DONE"
